# creature_team_怪物队伍表.xlsx — "new battle doing..." edit
#
# Transforms the old 5-column "top/mid/btm position" layout into the new
# 7-column 5-role battle-formation layout:
#   old: TeamID | desc | top#id_lv  | mid#id_lv  | btm#id_lv
#   new: TeamID | desc | l_pioneer#id_lv | r_pioneer#id_lv | commander#id_lv | l_guarder#id_lv | r_guarder#id_lv
# (displayed headers: 上位/中位/下位 -> 左先锋/右先锋/主将/左辅助/右辅助)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Add two new columns (F, G) by cloning the format of the existing
#    C/D columns (so styles line up with the xfs already used by the
#    sheet: header cells use s=1/s=2 alternating, the "map" row uses
#    s=3, and so do all the data rows).
# ---------------------------------------------------------------------
$ws.Range("C1:C10").Copy() | Out-Null
$ws.Range("F1:F10").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("D1:D10").Copy() | Out-Null
$ws.Range("G1:G10").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# ---------------------------------------------------------------------
# 2) Row 1 — column headers (display names)
# ---------------------------------------------------------------------
$ws.Range("C1").Value = "左先锋"
$ws.Range("D1").Value = "右先锋"
$ws.Range("E1").Value = "主将"
$ws.Range("F1").Value = "左辅助"
$ws.Range("G1").Value = "右辅助"

# ---------------------------------------------------------------------
# 3) Row 2 — "jl" client/server marker row
# ---------------------------------------------------------------------
$ws.Range("F2").Value = "jl"
$ws.Range("G2").Value = "jl"

# ---------------------------------------------------------------------
# 4) Row 3 — field keys
# ---------------------------------------------------------------------
$ws.Range("C3").Value = "l_pioneer#id_lv"
$ws.Range("D3").Value = "r_pioneer#id_lv"
$ws.Range("E3").Value = "commander#id_lv"
$ws.Range("F3").Value = "l_guarder#id_lv"
$ws.Range("G3").Value = "r_guarder#id_lv"

# ---------------------------------------------------------------------
# 5) Row 4 — field types ("map" for all 5 role columns)
# ---------------------------------------------------------------------
$ws.Range("F4").Value = "map"
$ws.Range("G4").Value = "map"

# ---------------------------------------------------------------------
# 6) Data rows 5-10 — the actual creature-team formation values
# ---------------------------------------------------------------------
$ws.Range("C5").Value = "10001_1"
$ws.Range("D5").Value = "10002_1"
$ws.Range("E5").Value = "10003_1"
$ws.Range("F5").Value = "10004_1"
$ws.Range("G5").Value = "10005_1"

$ws.Range("C6").Value = "10005_1"
$ws.Range("D6").Value = "10004_1"
$ws.Range("E6").Value = "10003_1"
$ws.Range("F6").Value = "10002_1"
$ws.Range("G6").Value = "10001_1"

$ws.Range("D7").Value = "10003_1"
$ws.Range("E7").Value = "10003_1"
$ws.Range("G7").Value = "10003_1"

$ws.Range("D8").Value = "10004_1"
$ws.Range("E8").Value = "10004_1"
$ws.Range("G8").Value = "10004_1"

$ws.Range("C9").Value = "10005_1"
$ws.Range("D9").Value = "10005_1"
$ws.Range("E9").Value = "10005_1"
$ws.Range("F9").Value = "10005_1"
$ws.Range("G9").Value = "10005_1"

$ws.Range("C10").Value = "10006_1"
$ws.Range("D10").Value = "10006_1"
$ws.Range("E10").Value = "10006_1"
$ws.Range("F10").Value = "10006_1"
$ws.Range("G10").Value = "10006_1"

# ---------------------------------------------------------------------
# 7) Column widths — C:G now share one uniform best-fit width instead of
#    the old bespoke per-column widths (target raw OOXML width 17.25;
#    the Excel ColumnWidth property round-trips through a pixel quantum
#    of 1/7 character, so 16.4 is the input that lands closest to it).
# ---------------------------------------------------------------------
$ws.Range("C1:G10").ColumnWidth = 16.4

# ---------------------------------------------------------------------
# 8) Restore the view state the author left the sheet in.
# ---------------------------------------------------------------------
$ws.Range("G14").Select() | Out-Null
